$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content correction: C5 questionMark value 1 -> 0.25 ---
$ws.Range("C5").Value = 0.25

# --- Apply whole-number format to the isCorrect flag columns (G, I, K, M) ---
# Doing this column-by-column (not as one comma-joined multi-area range) so the
# engine rebuilds the style table the same way Excel does when each column is
# formatted individually.
$ws.Columns.Item(7).NumberFormat = "0"
$ws.Columns.Item(9).NumberFormat = "0"
$ws.Columns.Item(11).NumberFormat = "0"
$ws.Columns.Item(13).NumberFormat = "0"

# --- Row heights shrink (content/wrap re-measured after formatting pass) ---
$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 43.2

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- View: scroll back to top-left and move the selection to K3 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K3").Select()
